$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "TestValues"

# Row headers (shared strings) - order of first use matters for shared string table indices
$ws.Range("A2").Value = "Num_Iterations"
$ws.Range("A3").Value = "Num_dimesdions"
$ws.Range("A1").Value = "Num_Particles"

# Row 1: Num_Particles values 10..100
$ws.Range("B1").Value = 10
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 30
$ws.Range("E1").Value = 40
$ws.Range("F1").Value = 50
$ws.Range("G1").Value = 60
$ws.Range("H1").Value = 70
$ws.Range("I1").Value = 80
$ws.Range("J1").Value = 90
$ws.Range("K1").Value = 100

# Row 2: Num_Iterations values 100..1000
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = 300
$ws.Range("E2").Value = 400
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 600
$ws.Range("H2").Value = 700
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 900
$ws.Range("K2").Value = 1000

# Row 3: Num_dimesdions values 2..11
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 11

# Column widths (target stored widths: 16, 16.140625, 14.140625 - the engine
# quantizes ColumnWidth to 1/6-character pixel steps, so we pick the input
# that lands on the closest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 15.17
$ws.Columns.Item(2).ColumnWidth = 15.33
$ws.Columns.Item(3).ColumnWidth = 13.33

# Selection
$null = $ws.Range("B18").Select()
